$p = $ppt.ActivePresentation
$s16 = $p.Slides.Item(16)
$shp = $s16.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange
$chars = $tr.Characters(39, 42)
$chars.Text = "می‌خوایم که قابلیت‌‌های زیر رو داشته باشه:"
